$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_16_3_0"
$ws.Range("B2").Value = [double]"0.9999911475978447"
$ws.Range("C2").Value = [double]"0.9991050381053288"
$ws.Range("D2").Value = [double]"0.9999626817613565"
$ws.Range("E2").Value = [double]"0.9999969221167055"
$ws.Range("F2").Value = [double]"0.999980100715405"
$ws.Range("G2").Value = [double]"8.263328836443779e-06"
$ws.Range("H2").Value = [double]"0.0008354076443886224"
$ws.Range("I2").Value = [double]"3.806218146854155e-05"
$ws.Range("J2").Value = [double]"3.105863851827505e-06"
$ws.Range("K2").Value = [double]"2.058402266018452e-05"
$ws.Range("L2").Value = [double]"0.0002264236711867637"
$ws.Range("M2").Value = [double]"0.002874600639470426"
$ws.Range("N2").Value = [double]"1.000023606405748"
$ws.Range("O2").Value = [double]"0.002996978339936341"
$ws.Range("P2").Value = [double]"89.40736608956689"
$ws.Range("Q2").Value = [double]"129.6302683102175"

$ws.Range("A3").Value = "model_16_3_1"
$ws.Range("B3").Value = [double]"0.9999920041558984"
$ws.Range("C3").Value = [double]"0.99910366263786"
$ws.Range("D3").Value = [double]"0.9999659021094975"
$ws.Range("E3").Value = [double]"0.9999966874111882"
$ws.Range("F3").Value = [double]"0.9999815738770962"
$ws.Range("G3").Value = [double]"7.463769491874863e-06"
$ws.Range("H3").Value = [double]"0.0008366915828947142"
$ws.Range("I3").Value = [double]"3.477763536478169e-05"
$ws.Range("J3").Value = [double]"3.342703040514164e-06"
$ws.Range("K3").Value = [double]"1.906016920264792e-05"
$ws.Range("L3").Value = [double]"0.0002167208326791978"
$ws.Range("M3").Value = [double]"0.002731990024117011"
$ws.Range("N3").Value = [double]"1.000021322250938"
$ws.Range("O3").Value = [double]"0.00284829649544266"
$ws.Range("P3").Value = [double]"89.61089995509906"
$ws.Range("Q3").Value = [double]"129.8338021757497"

$ws.Range("A4").Value = "model_16_3_2"
$ws.Range("B4").Value = [double]"0.9999926992670658"
$ws.Range("C4").Value = [double]"0.9991023814842022"
$ws.Range("D4").Value = [double]"0.999968532325481"
$ws.Range("E4").Value = [double]"0.999996476472498"
$ws.Range("F4").Value = [double]"0.9999827676942769"
$ws.Range("G4").Value = [double]"6.814913728929273e-06"
$ws.Range("H4").Value = [double]"0.000837887483598195"
$ws.Range("I4").Value = [double]"3.209498576220556e-05"
$ws.Range("J4").Value = [double]"3.555559341375273e-06"
$ws.Range("K4").Value = [double]"1.782527255179042e-05"
$ws.Range("L4").Value = [double]"0.0002078580626749844"
$ws.Range("M4").Value = [double]"0.002610538972880749"
$ws.Range("N4").Value = [double]"1.000019468621158"
$ws.Range("O4").Value = [double]"0.002721675021516935"
$ws.Range("P4").Value = [double]"89.79279430378762"
$ws.Range("Q4").Value = [double]"130.0156965244383"

$ws.Range("A5").Value = "model_16_3_3"
$ws.Range("B5").Value = [double]"0.999993266335392"
$ws.Range("C5").Value = [double]"0.9991011804997167"
$ws.Range("D5").Value = [double]"0.9999706920662961"
$ws.Range("E5").Value = [double]"0.9999962892491258"
$ws.Range("F5").Value = [double]"0.9999837411331036"
$ws.Range("G5").Value = [double]"6.285580337916873e-06"
$ws.Range("H5").Value = [double]"0.0008390085498982103"
$ws.Range("I5").Value = [double]"2.989219029768033e-05"
$ws.Range("J5").Value = [double]"3.744484732123362e-06"
$ws.Range("K5").Value = [double]"1.681833751490184e-05"
$ws.Range("L5").Value = [double]"0.0001997172082916236"
$ws.Range("M5").Value = [double]"0.002507105968625354"
$ws.Range("N5").Value = [double]"1.000017956438955"
$ws.Range("O5").Value = [double]"0.002613838660134552"
$ws.Range("P5").Value = [double]"89.95450476641049"
$ws.Range("Q5").Value = [double]"130.1774069870611"

$ws.Range("A6").Value = "model_16_3_4"
$ws.Range("B6").Value = [double]"0.9999937327263918"
$ws.Range("C6").Value = [double]"0.9991000453595962"
$ws.Range("D6").Value = [double]"0.9999725019440968"
$ws.Range("E6").Value = [double]"0.9999961117167364"
$ws.Range("F6").Value = [double]"0.9999845468152264"
$ws.Range("G6").Value = [double]"5.850224811734915e-06"
$ws.Range("H6").Value = [double]"0.000840068153373838"
$ws.Range("I6").Value = [double]"2.804623240166881e-05"
$ws.Range("J6").Value = [double]"3.923631040795344e-06"
$ws.Range("K6").Value = [double]"1.598493172123208e-05"
$ws.Range("L6").Value = [double]"0.0001922977285361094"
$ws.Range("M6").Value = [double]"0.002418723798149535"
$ws.Range("N6").Value = [double]"1.000016712729622"
$ws.Range("O6").Value = [double]"0.002521693877685264"
$ws.Range("P6").Value = [double]"90.098060936205"
$ws.Range("Q6").Value = [double]"130.3209631568556"

$ws.Range("A7").Value = "model_16_3_5"
$ws.Range("B7").Value = [double]"0.999994111005992"
$ws.Range("C7").Value = [double]"0.9990989966845555"
$ws.Range("D7").Value = [double]"0.9999739722438663"
$ws.Range("E7").Value = [double]"0.9999959614895927"
$ws.Range("F7").Value = [double]"0.9999851984027334"
$ws.Range("G7").Value = [double]"5.497117409485894e-06"
$ws.Range("H7").Value = [double]"0.0008410470454927744"
$ws.Range("I7").Value = [double]"2.654662205911985e-05"
$ws.Range("J7").Value = [double]"4.075223876032954e-06"
$ws.Range("K7").Value = [double]"1.53109229675764e-05"
$ws.Range("L7").Value = [double]"0.000185472179555587"
$ws.Range("M7").Value = [double]"0.002344593229002825"
$ws.Range("N7").Value = [double]"1.000015703984021"
$ws.Range("O7").Value = [double]"0.002444407416738545"
$ws.Range("P7").Value = [double]"90.22257342096864"
$ws.Range("Q7").Value = [double]"130.4454756416193"

$ws.Range("A8").Value = "model_16_3_6"
$ws.Range("B8").Value = [double]"0.9999944210319018"
$ws.Range("C8").Value = [double]"0.9990980212704257"
$ws.Range("D8").Value = [double]"0.9999751984112603"
$ws.Range("E8").Value = [double]"0.9999958209727172"
$ws.Range("F8").Value = [double]"0.9999857343684881"
$ws.Range("G8").Value = [double]"5.20772183128143e-06"
$ws.Range("H8").Value = [double]"0.0008419575517672926"
$ws.Range("I8").Value = [double]"2.529601089529365e-05"
$ws.Range("J8").Value = [double]"4.217018168571559e-06"
$ws.Range("K8").Value = [double]"1.475651453193261e-05"
$ws.Range("L8").Value = [double]"0.0001792711707536339"
$ws.Range("M8").Value = [double]"0.002282043345618446"
$ws.Range("N8").Value = [double]"1.000014877248262"
$ws.Range("O8").Value = [double]"0.002379194655322384"
$ws.Range("P8").Value = [double]"90.33073613259359"
$ws.Range("Q8").Value = [double]"130.5536383532442"

$ws.Range("A9").Value = "model_16_3_7"
$ws.Range("B9").Value = [double]"0.9999946721794911"
$ws.Range("C9").Value = [double]"0.9990971335331363"
$ws.Range("D9").Value = [double]"0.9999761940896789"
$ws.Range("E9").Value = [double]"0.9999957023926674"
$ws.Range("F9").Value = [double]"0.9999861674023746"
$ws.Range("G9").Value = [double]"4.973286580784456e-06"
$ws.Range("H9").Value = [double]"0.0008427862155597256"
$ws.Range("I9").Value = [double]"2.428048352764854e-05"
$ws.Range("J9").Value = [double]"4.336676211053874e-06"
$ws.Range("K9").Value = [double]"1.43085798693512e-05"
$ws.Range("L9").Value = [double]"0.0001736213406782131"
$ws.Range("M9").Value = [double]"0.002230086675621478"
$ws.Range("N9").Value = [double]"1.000014207521357"
$ws.Range("O9").Value = [double]"0.002325026082318511"
$ws.Range("P9").Value = [double]"90.42285930509502"
$ws.Range("Q9").Value = [double]"130.6457615257457"

$ws.Range("A10").Value = "model_16_3_8"
$ws.Range("B10").Value = [double]"0.9999948793907101"
$ws.Range("C10").Value = [double]"0.9990963058275781"
$ws.Range("D10").Value = [double]"0.9999770383034909"
$ws.Range("E10").Value = [double]"0.9999955892544533"
$ws.Range("F10").Value = [double]"0.9999865284180401"
$ws.Range("G10").Value = [double]"4.779864003359813e-06"
$ws.Range("H10").Value = [double]"0.0008435588423661621"
$ws.Range("I10").Value = [double]"2.341944022876127e-05"
$ws.Range("J10").Value = [double]"4.450842947057679e-06"
$ws.Range("K10").Value = [double]"1.393514158790948e-05"
$ws.Range("L10").Value = [double]"0.000168506426149669"
$ws.Range("M10").Value = [double]"0.002186290008978638"
$ws.Range("N10").Value = [double]"1.000013654958106"
$ws.Range("O10").Value = [double]"0.002279364900905087"
$ws.Range("P10").Value = [double]"90.50219692609122"
$ws.Range("Q10").Value = [double]"130.7250991467419"

$ws.Range("A11").Value = "model_16_3_9"
$ws.Range("B11").Value = [double]"0.9999950474579137"
$ws.Range("C11").Value = [double]"0.9990955463978271"
$ws.Range("D11").Value = [double]"0.9999777279827216"
$ws.Range("E11").Value = [double]"0.999995489587635"
$ws.Range("F11").Value = [double]"0.9999868198184875"
$ws.Range("G11").Value = [double]"4.622980646075969e-06"
$ws.Range("H11").Value = [double]"0.0008442677367037043"
$ws.Range("I11").Value = [double]"2.271601217360586e-05"
$ws.Range("J11").Value = [double]"4.551415820917982e-06"
$ws.Range("K11").Value = [double]"1.363371399726192e-05"
$ws.Range("L11").Value = [double]"0.0001638477480463509"
$ws.Range("M11").Value = [double]"0.002150111775251689"
$ws.Range("N11").Value = [double]"1.000013206778897"
$ws.Range("O11").Value = [double]"0.002241646484869114"
$ws.Range("P11").Value = [double]"90.56894179881905"
$ws.Range("Q11").Value = [double]"130.7918440194697"

$ws.Range("A12").Value = "model_16_3_10"
$ws.Range("B12").Value = [double]"0.9999951859047793"
$ws.Range("C12").Value = [double]"0.9990948539181025"
$ws.Range("D12").Value = [double]"0.9999783092627047"
$ws.Range("E12").Value = [double]"0.999995398760758"
$ws.Range("F12").Value = [double]"0.9999870620895398"
$ws.Range("G12").Value = [double]"4.493746574185288e-06"
$ws.Range("H12").Value = [double]"0.0008449141361302653"
$ws.Range("I12").Value = [double]"2.212314431570388e-05"
$ws.Range("J12").Value = [double]"4.643068390855881e-06"
$ws.Range("K12").Value = [double]"1.338310635327988e-05"
$ws.Range("L12").Value = [double]"0.0001596163218539301"
$ws.Range("M12").Value = [double]"0.002119845884536253"
$ws.Range("N12").Value = [double]"1.000012837587255"
$ws.Range("O12").Value = [double]"0.00221009211252689"
$ws.Range("P12").Value = [double]"90.62564755565666"
$ws.Range("Q12").Value = [double]"130.8485497763073"

$ws.Range("A13").Value = "model_16_3_11"
$ws.Range("B13").Value = [double]"0.9999952980562711"
$ws.Range("C13").Value = [double]"0.9990942245706343"
$ws.Range("D13").Value = [double]"0.9999787823074177"
$ws.Range("E13").Value = [double]"0.9999953242340487"
$ws.Range("F13").Value = [double]"0.9999872589508489"
$ws.Range("G13").Value = [double]"4.389058079450089e-06"
$ws.Range("H13").Value = [double]"0.0008455016043666111"
$ws.Range("I13").Value = [double]"2.16406694090597e-05"
$ws.Range("J13").Value = [double]"4.71827261085943e-06"
$ws.Range("K13").Value = [double]"1.317947100995957e-05"
$ws.Range("L13").Value = [double]"0.0001557656107382681"
$ws.Range("M13").Value = [double]"0.002095007894841947"
$ws.Range("N13").Value = [double]"1.000012538516611"
$ws.Range("O13").Value = [double]"0.002184196718189572"
$ws.Range("P13").Value = [double]"90.67279182877664"
$ws.Range("Q13").Value = [double]"130.8956940494272"

$ws.Range("A14").Value = "model_16_3_12"
$ws.Range("B14").Value = [double]"0.9999953904942942"
$ws.Range("C14").Value = [double]"0.9990936469621484"
$ws.Range("D14").Value = [double]"0.9999791826523881"
$ws.Range("E14").Value = [double]"0.9999952528355832"
$ws.Range("F14").Value = [double]"0.999987421496788"
$ws.Range("G14").Value = [double]"4.302771242497516e-06"
$ws.Range("H14").Value = [double]"0.000846040776534165"
$ws.Range("I14").Value = [double]"2.123234352160186e-05"
$ws.Range("J14").Value = [double]"4.790320148667976e-06"
$ws.Range("K14").Value = [double]"1.301133183513492e-05"
$ws.Range("L14").Value = [double]"0.0001522707876765186"
$ws.Range("M14").Value = [double]"0.002074312233608411"
$ws.Range("N14").Value = [double]"1.000012292015215"
$ws.Range("O14").Value = [double]"0.002162619999811399"
$ws.Range("P14").Value = [double]"90.71250253570032"
$ws.Range("Q14").Value = [double]"130.9354047563509"

$ws.Range("A15").Value = "model_16_3_13"
$ws.Range("B15").Value = [double]"0.9999954656099161"
$ws.Range("C15").Value = [double]"0.9990931222527707"
$ws.Range("D15").Value = [double]"0.9999795160214361"
$ws.Range("E15").Value = [double]"0.9999951910834532"
$ws.Range("F15").Value = [double]"0.9999875557284924"
$ws.Range("G15").Value = [double]"4.232654106594881e-06"
$ws.Range("H15").Value = [double]"0.0008465305697062136"
$ws.Range("I15").Value = [double]"2.089232872820522e-05"
$ws.Range("J15").Value = [double]"4.85263365767333e-06"
$ws.Range("K15").Value = [double]"1.287248119293928e-05"
$ws.Range("L15").Value = [double]"0.0001491357862309553"
$ws.Range("M15").Value = [double]"0.002057341514332242"
$ws.Range("N15").Value = [double]"1.00001209170689"
$ws.Range("O15").Value = [double]"0.002144926801881411"
$ws.Range("P15").Value = [double]"90.74536262660385"
$ws.Range("Q15").Value = [double]"130.9682648472545"

$ws.Range("A16").Value = "model_16_3_14"
$ws.Range("B16").Value = [double]"0.9999955276857565"
$ws.Range("C16").Value = [double]"0.9990926478603426"
$ws.Range("D16").Value = [double]"0.999979789423073"
$ws.Range("E16").Value = [double]"0.999995137795345"
$ws.Range("F16").Value = [double]"0.9999876645244937"
$ws.Range("G16").Value = [double]"4.174709034403751e-06"
$ws.Range("H16").Value = [double]"0.0008469733942144238"
$ws.Range("I16").Value = [double]"2.061347680226235e-05"
$ws.Range("J16").Value = [double]"4.906406199805502e-06"
$ws.Range("K16").Value = [double]"1.275994150103393e-05"
$ws.Range("L16").Value = [double]"0.0001462728454265286"
$ws.Range("M16").Value = [double]"0.002043210472370321"
$ws.Range("N16").Value = [double]"1.000011926171316"
$ws.Range("O16").Value = [double]"0.002130194172207882"
$ws.Range("P16").Value = [double]"90.77293178889693"
$ws.Range("Q16").Value = [double]"130.9958340095475"

$ws.Range("A17").Value = "model_16_3_15"
$ws.Range("B17").Value = [double]"0.9999955780147405"
$ws.Range("C17").Value = [double]"0.9990922132195543"
$ws.Range("D17").Value = [double]"0.9999800223506805"
$ws.Range("E17").Value = [double]"0.9999950871780005"
$ws.Range("F17").Value = [double]"0.999987754669364"
$ws.Range("G17").Value = [double]"4.12772913708857e-06"
$ws.Range("H17").Value = [double]"0.0008473791123116702"
$ws.Range("I17").Value = [double]"2.037590575957777e-05"
$ws.Range("J17").Value = [double]"4.957483698675872e-06"
$ws.Range("K17").Value = [double]"1.266669472912682e-05"
$ws.Range("L17").Value = [double]"0.0001436885624907152"
$ws.Range("M17").Value = [double]"0.002031681357174045"
$ws.Range("N17").Value = [double]"1.000011791960692"
$ws.Range("O17").Value = [double]"0.002118174238709142"
$ws.Range("P17").Value = [double]"90.79556629590641"
$ws.Range("Q17").Value = [double]"131.018468516557"

$ws.Range("A18").Value = "model_16_3_16"
$ws.Range("B18").Value = [double]"0.9999956200599946"
$ws.Range("C18").Value = [double]"0.9990918245270682"
$ws.Range("D18").Value = [double]"0.9999802191272872"
$ws.Range("E18").Value = [double]"0.9999950425283021"
$ws.Range("F18").Value = [double]"0.9999878299024623"
$ws.Range("G18").Value = [double]"4.088481737896318e-06"
$ws.Range("H18").Value = [double]"0.000847741939685833"
$ws.Range("I18").Value = [double]"2.017520639152423e-05"
$ws.Range("J18").Value = [double]"5.002539300605711e-06"
$ws.Range("K18").Value = [double]"1.258887284606497e-05"
$ws.Range("L18").Value = [double]"0.0001413560819110566"
$ws.Range("M18").Value = [double]"0.002021999440627103"
$ws.Range("N18").Value = [double]"1.000011679840014"
$ws.Range("O18").Value = [double]"0.002108080142930467"
$ws.Range("P18").Value = [double]"90.81467374010907"
$ws.Range("Q18").Value = [double]"131.0375759607597"

$ws.Range("A19").Value = "model_16_3_17"
$ws.Range("B19").Value = [double]"0.9999956536380694"
$ws.Range("C19").Value = [double]"0.9990914693149036"
$ws.Range("D19").Value = [double]"0.9999803829454436"
$ws.Range("E19").Value = [double]"0.9999950044923746"
$ws.Range("F19").Value = [double]"0.999987892112881"
$ws.Range("G19").Value = [double]"4.057138078928856e-06"
$ws.Range("H19").Value = [double]"0.0008480735146495981"
$ws.Range("I19").Value = [double]"2.000812250382407e-05"
$ws.Range("J19").Value = [double]"5.040921006798866e-06"
$ws.Range("K19").Value = [double]"1.252452175531147e-05"
$ws.Range("L19").Value = [double]"0.0001392450551657073"
$ws.Range("M19").Value = [double]"0.00201423386897571"
$ws.Range("N19").Value = [double]"1.000011590298482"
$ws.Range("O19").Value = [double]"0.002099983974816924"
$ws.Range("P19").Value = [double]"90.83006547916385"
$ws.Range("Q19").Value = [double]"131.0529676998145"

$ws.Range("A20").Value = "model_16_3_18"
$ws.Range("B20").Value = [double]"0.9999956818648695"
$ws.Range("C20").Value = [double]"0.9990911500662597"
$ws.Range("D20").Value = [double]"0.9999805192873141"
$ws.Range("E20").Value = [double]"0.9999949703446922"
$ws.Range("F20").Value = [double]"0.9999879426739385"
$ws.Range("G20").Value = [double]"4.030789600053497e-06"
$ws.Range("H20").Value = [double]"0.0008483715192452071"
$ws.Range("I20").Value = [double]"1.986906264438601e-05"
$ws.Range("J20").Value = [double]"5.075379120350778e-06"
$ws.Range("K20").Value = [double]"1.247222088236839e-05"
$ws.Range("L20").Value = [double]"0.0001373420390465239"
$ws.Range("M20").Value = [double]"0.002007682644257677"
$ws.Range("N20").Value = [double]"1.000011515027015"
$ws.Range("O20").Value = [double]"0.002093153850899739"
$ws.Range("P20").Value = [double]"90.84309654132639"
$ws.Range("Q20").Value = [double]"131.065998761977"

$ws.Range("A21").Value = "model_16_3_19"
$ws.Range("B21").Value = [double]"0.9999957052443977"
$ws.Range("C21").Value = [double]"0.9990908581840229"
$ws.Range("D21").Value = [double]"0.9999806385101675"
$ws.Range("E21").Value = [double]"0.9999949391931053"
$ws.Range("F21").Value = [double]"0.9999879862566454"
$ws.Range("G21").Value = [double]"4.008965836751301e-06"
$ws.Range("H21").Value = [double]"0.0008486439785011085"
$ws.Range("I21").Value = [double]"1.974746307145237e-05"
$ws.Range("J21").Value = [double]"5.106813901491782e-06"
$ws.Range("K21").Value = [double]"1.242713848647208e-05"
$ws.Range("L21").Value = [double]"0.0001356237141184489"
$ws.Range("M21").Value = [double]"0.002002240204558709"
$ws.Range("N21").Value = [double]"1.000011452681606"
$ws.Range("O21").Value = [double]"0.002087479715275384"
$ws.Range("P21").Value = [double]"90.85395449195741"
$ws.Range("Q21").Value = [double]"131.076856712608"

$ws.Range("A22").Value = "model_16_3_20"
$ws.Range("B22").Value = [double]"0.9999957242476973"
$ws.Range("C22").Value = [double]"0.9990905965806868"
$ws.Range("D22").Value = [double]"0.9999807378398824"
$ws.Range("E22").Value = [double]"0.9999949101248259"
$ws.Range("F22").Value = [double]"0.9999880210481281"
$ws.Range("G22").Value = [double]"3.991227090660813e-06"
$ws.Range("H22").Value = [double]"0.0008488881737322692"
$ws.Range("I22").Value = [double]"1.964615320883204e-05"
$ws.Range("J22").Value = [double]"5.136146436104547e-06"
$ws.Range("K22").Value = [double]"1.23911498224683e-05"
$ws.Range("L22").Value = [double]"0.0001340693071656698"
$ws.Range("M22").Value = [double]"0.001997805568783112"
$ws.Range("N22").Value = [double]"1.00001140200614"
$ws.Range("O22").Value = [double]"0.002082856287873858"
$ws.Range("P22").Value = [double]"90.8628236656494"
$ws.Range("Q22").Value = [double]"131.0857258863"

$ws.Range("A23").Value = "model_16_3_21"
$ws.Range("B23").Value = [double]"0.9999957403412444"
$ws.Range("C23").Value = [double]"0.9990903606231543"
$ws.Range("D23").Value = [double]"0.9999808197895046"
$ws.Range("E23").Value = [double]"0.9999948884149666"
$ws.Range("F23").Value = [double]"0.9999880508603229"
$ws.Range("G23").Value = [double]"3.976204470956059e-06"
$ws.Range("H23").Value = [double]"0.0008491084297316244"
$ws.Range("I23").Value = [double]"1.956256991266106e-05"
$ws.Range("J23").Value = [double]"5.158053656340314e-06"
$ws.Range("K23").Value = [double]"1.236031178450069e-05"
$ws.Range("L23").Value = [double]"0.0001326642424330314"
$ws.Range("M23").Value = [double]"0.00199404224402495"
$ws.Range("N23").Value = [double]"1.000011359090015"
$ws.Range("O23").Value = [double]"0.002078932750589585"
$ws.Range("P23").Value = [double]"90.87036568839088"
$ws.Range("Q23").Value = [double]"131.0932679090415"

$ws.Range("A24").Value = "model_16_3_22"
$ws.Range("B24").Value = [double]"0.9999957533316737"
$ws.Range("C24").Value = [double]"0.9990901456242437"
$ws.Range("D24").Value = [double]"0.9999808948876819"
$ws.Range("E24").Value = [double]"0.9999948642615037"
$ws.Range("F24").Value = [double]"0.9999880761028377"
$ws.Range("G24").Value = [double]"3.96407847543067e-06"
$ws.Range("H24").Value = [double]"0.0008493091217772422"
$ws.Range("I24").Value = [double]"1.948597464563328e-05"
$ws.Range("J24").Value = [double]"5.182426694588202e-06"
$ws.Range("K24").Value = [double]"1.233420067011074e-05"
$ws.Range("L24").Value = [double]"0.0001314014010776326"
$ws.Range("M24").Value = [double]"0.001990999366004588"
$ws.Range("N24").Value = [double]"1.00001132444887"
$ws.Range("O24").Value = [double]"0.00207576033095227"
$ws.Range("P24").Value = [double]"90.87647428932308"
$ws.Range("Q24").Value = [double]"131.0993765099737"

$ws.Range("A25").Value = "model_16_3_23"
$ws.Range("B25").Value = [double]"0.9999957647030231"
$ws.Range("C25").Value = [double]"0.9990899530084837"
$ws.Range("D25").Value = [double]"0.9999809584037614"
$ws.Range("E25").Value = [double]"0.9999948441936329"
$ws.Range("F25").Value = [double]"0.9999880976281353"
$ws.Range("G25").Value = [double]"3.953463819858684e-06"
$ws.Range("H25").Value = [double]"0.0008494889201343099"
$ws.Range("I25").Value = [double]"1.942119236673043e-05"
$ws.Range("J25").Value = [double]"5.202676999293394e-06"
$ws.Range("K25").Value = [double]"1.231193468301191e-05"
$ws.Range("L25").Value = [double]"0.0001302546443771097"
$ws.Range("M25").Value = [double]"0.00198833191893574"
$ws.Range("N25").Value = [double]"1.000011294125272"
$ws.Range("O25").Value = [double]"0.002072979325139325"
$ws.Range("P25").Value = [double]"90.88183689379579"
$ws.Range("Q25").Value = [double]"131.1047391144464"

$ws.Range("A26").Value = "model_16_3_24"
$ws.Range("B26").Value = [double]"0.9999957732449217"
$ws.Range("C26").Value = [double]"0.9990897798337718"
$ws.Range("D26").Value = [double]"0.9999810077181477"
$ws.Range("E26").Value = [double]"0.9999948308699746"
$ws.Range("F26").Value = [double]"0.9999881154415194"
$ws.Range("G26").Value = [double]"3.94549033251495e-06"
$ws.Range("H26").Value = [double]"0.0008496505711264018"
$ws.Range("I26").Value = [double]"1.937089489324468e-05"
$ws.Range("J26").Value = [double]"5.216121780841699e-06"
$ws.Range("K26").Value = [double]"1.229350833704319e-05"
$ws.Range("L26").Value = [double]"0.0001292388613048028"
$ws.Range("M26").Value = [double]"0.001986325837448366"
$ws.Range("N26").Value = [double]"1.000011271346875"
$ws.Range("O26").Value = [double]"0.002070887840609871"
$ws.Range("P26").Value = [double]"90.88587463848864"
$ws.Range("Q26").Value = [double]"131.1087768591393"
